$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Restructured control processing -> container: append the new
# "Abwesenheiten" (Normalfall Anlage) test case as row 3, following the
# same layout as the existing rows.
$ws.Range("A3").Value = "003_Profil_003_Abwesenheiten_Normalfall_Anlage"
$ws.Range("B3").Value = "var003_Profil_003_Abwesenheiten_Normalfall_Anlage"
$ws.Range("C3").Value = "001_Login_001_Successful"
$ws.Range("D3").Value = "Menueauswahl Mein Profil"
$ws.Range("E3").Value = "003_Profil_003_Abwesenheiten_Normalfall_Anlage"
$ws.Range("F3").Value = "Logoff"

$ws.Range("E3").Select()
